# Weekly fruit/vegetable update: two new price records (rows 396-397) are
# inserted into the "Pimiento" (Macroferia Regional de Talca) sheet, pushing
# all existing records from row 396 onward down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 396, shifting old rows 396:420 down to 398:422.
$ws.Rows("396:397").Insert()

# --- New row 396 --------------------------------------------------------
$ws.Cells.Item(396, 1).Value  = 5
$ws.Cells.Item(396, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(396, 3).Value  = "Maule"
$ws.Cells.Item(396, 4).Value  = 44585
$ws.Cells.Item(396, 5).Value  = 7
$ws.Cells.Item(396, 6).Value  = 100112002
$ws.Cells.Item(396, 7).Value  = "Pimiento"
$ws.Cells.Item(396, 8).Value  = "Cuatro cascos verde"
$ws.Cells.Item(396, 9).Value  = "Primera"
$ws.Cells.Item(396, 10).Value = 300
$ws.Cells.Item(396, 11).Value = 7000
$ws.Cells.Item(396, 12).Value = 7000
$ws.Cells.Item(396, 13).Value = 7000
$ws.Cells.Item(396, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(396, 15).Value = "Región del Maule"
$ws.Cells.Item(396, 16).Value = 467
$ws.Cells.Item(396, 17).Value = 15
$ws.Cells.Item(396, 18).Value = "Hortaliza"

# --- New row 397 --------------------------------------------------------
$ws.Cells.Item(397, 1).Value  = 5
$ws.Cells.Item(397, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(397, 3).Value  = "Maule"
$ws.Cells.Item(397, 4).Value  = 44585
$ws.Cells.Item(397, 5).Value  = 7
$ws.Cells.Item(397, 6).Value  = 100112002
$ws.Cells.Item(397, 7).Value  = "Pimiento"
$ws.Cells.Item(397, 8).Value  = "Morrón rojo"
$ws.Cells.Item(397, 9).Value  = "Primera"
$ws.Cells.Item(397, 10).Value = 150
$ws.Cells.Item(397, 11).Value = 15000
$ws.Cells.Item(397, 12).Value = 15000
$ws.Cells.Item(397, 13).Value = 15000
$ws.Cells.Item(397, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(397, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(397, 16).Value = 750
$ws.Cells.Item(397, 17).Value = 20
$ws.Cells.Item(397, 18).Value = "Hortaliza"
